# Update market_health_data.xlsx per commit "Update data: 2025-11-05 09:50"
$wb = $excel.ActiveWorkbook

# 1) Update the "Last Updated" timestamp on the Metadata sheet
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "05 Nov 2025, 09:49 AM"

# 2) Update the "1 Year" column (F) values on the Industry Analysis sheet
$wsIndustry = $wb.Worksheets.Item("Industry Analysis")

$newValues = @(
    21.3,
    -4.3927,
    35.9445,
    -51.0482,
    57.2275,
    -9.640700000000001,
    -6.1449,
    36.9733,
    -4.7026,
    46.5317,
    -2.102,
    17.4681,
    -33.0245,
    1.0205,
    2.0426,
    -16.2411,
    7.4627,
    -25.798,
    47.7485,
    19.5587,
    76.5603,
    -54.2675,
    -0.8811,
    4.8518,
    3.6831,
    -34.0874,
    -11.9893,
    -12.994,
    25.5415,
    56.5088,
    2.0908,
    -4.7193,
    22.8807,
    5.3359,
    -5.1995,
    -5.6238,
    -22.595,
    10.8405,
    -7.5963,
    -4.552,
    22.3098,
    14.0694,
    -9.6066,
    27.639,
    -6.3484,
    -40.5302,
    -29.7988,
    -24.0791,
    -49.1803,
    -51.6023,
    -34.4756,
    -11.5478,
    -2.3796,
    -15.4382,
    -27.6987,
    -27.1559,
    -2.1585,
    -23.0964,
    -13.3217,
    -8.1496,
    -16.0695,
    -12.5465,
    47.7264,
    -42.4232,
    11.3291,
    14.3746,
    32.6702,
    -17.0097,
    -13.5162,
    11.4259,
    2.6754,
    -11.1574,
    -13.2502,
    24.7078,
    53.3554
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $wsIndustry.Cells.Item($row, 6).Value = $newValues[$i]
}
